$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store plain text values (e.g. "3.563.01",
# "  +0.40%  ") that must not be auto-converted to numbers/dates by Excel's
# COM value coercion, so force the number format to Text before writing.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '95.856.14'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.563.01'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '238.82'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '654.07'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  +9.54%  '
$ws.Range('D8').Value = '0.403'
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  +4.73%  '
$ws.Range('D11').Value = '3.561.07'
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '43.27'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').Value = '6.39'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '4.223.81'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').Value = '95.797.12'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '3.566.63'
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('D19').Value = '7.74'
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('D20').Value = '12.60'
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('D21').Value = '17.66'
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('E22').Value = '  +5.83%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '501.23'
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '3.37'
$ws.Range('E24').Value = '  -7.25%  '
$ws.Range('D25').Value = '6.93'
$ws.Range('E25').Value = '  +5.21%  '
$ws.Range('D26').Value = '0.0000196'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('D27').Value = '95.56'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('E28').Value = '  +0.24%  '
$ws.Range('D29').Value = '3.755.31'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.152'
$ws.Range('E30').Value = '  +9.63%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.99'
$ws.Range('E31').Value = '  -5.00%  '
$ws.Range('D32').Value = '11.31'
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('D36').Value = '31.23'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').Value = '8.77'
$ws.Range('E37').Value = '  +7.37%  '
$ws.Range('D38').Value = '612.14'
$ws.Range('E38').Value = '  +6.98%  '
$ws.Range('D39').Value = '0.563'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').Value = '1.62'
$ws.Range('E40').Value = '  +9.90%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = '0.901'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('E44').Value = '  +5.19%  '
$ws.Range('D45').Value = '5.67'
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').Value = '23.51'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0419'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '33.87'
$ws.Range('E48').Value = '  -3.39%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.26'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '8.14'
$ws.Range('E51').Value = '  +0.94%  '
